# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The 'K' column (column G) holds recomputed strike-count ("s_vals") figures
# for each trade row. Replace the previously-written values with the
# freshly calculated ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 3
    16 = 3
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 3
    34 = 2
    35 = 2
    36 = 2
    37 = 1
    38 = 2
    39 = 0
    40 = 2
    41 = 2
    42 = 1
    43 = 2
    44 = 2
    45 = 2
    46 = 1
    47 = 3
    48 = 0
    49 = 1
    50 = 3
    51 = 0
    52 = 0
    53 = 0
    54 = 1
    55 = 0
    56 = 0
    57 = 1
    58 = 0
    59 = 0
    60 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}

Write-Host "Updated K column (G2:G60) with recalculated s_vals"
